$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.008
$ws.Range("C3").Value = -12.314
$ws.Range("E8").Value = 16.597
$ws.Range("E11").Value = 16.41
$ws.Range("A12").Value = -21.654
$ws.Range("B14").Value = 5.598999999999999
$ws.Range("E14").Value = 16.814
$ws.Range("E15").Value = 16.177
$ws.Range("B26").Value = 5.95
$ws.Range("C30").Value = -12.941
$ws.Range("B31").Value = 6.167999999999999
$ws.Range("A32").Value = -21.441
$ws.Range("B35").Value = 8.095000000000001
$ws.Range("A36").Value = -21.12
$ws.Range("E36").Value = 16.457
$ws.Range("B37").Value = 8.260000000000002
$ws.Range("A38").Value = -20.363
$ws.Range("C44").Value = -12.289
$ws.Range("B45").Value = 5.712000000000001
$ws.Range("A46").Value = -21.591
$ws.Range("A54").Value = -21.876
$ws.Range("A55").Value = -21.921
$ws.Range("B57").Value = 6.090000000000001
$ws.Range("C58").Value = -12.609
$ws.Range("E64").Value = 17.254
$ws.Range("A67").Value = -21.6
$ws.Range("A69").Value = -21.651
$ws.Range("A72").Value = -21.55
$ws.Range("C84").Value = -13.659
$ws.Range("C89").Value = -10.717
$ws.Range("E89").Value = 16.997
$ws.Range("A91").Value = -21.761
$ws.Range("C91").Value = -11.186
$ws.Range("C92").Value = -11.603
$ws.Range("A99").Value = -20.559
$ws.Range("B100").Value = 6.165000000000001
$ws.Range("B102").Value = 7.334999999999999
$ws.Range("C102").Value = -12.808
